$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the "Trial 2 / C++" block (old row 9),
# which pushes that block (old rows 9-14) down by one (new rows 10-15).
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with standard-deviation stats for the
# first ("Trial 1 / C") block (rows 2-6).
$ws.Range("C8").Value = "stdev"
$ws.Range("D8").Formula = "=_xlfn.STDEV.S(D2:D6)"
$ws.Range("E8").Formula = "=_xlfn.STDEV.S(E2:E6)"
$ws.Range("F8").Formula = "=_xlfn.STDEV.S(F2:F6)"
$ws.Range("G8").Formula = "=_xlfn.STDEV.S(G2:G6)"

# Append a new row 16 with standard-deviation stats for the second
# ("Trial 2 / C++") block, now living at rows 10-14 (averages at row 15).
$ws.Range("C16").Value = "stdev"
$ws.Range("D16").Formula = "=_xlfn.STDEV.S(D10:D14)"
$ws.Range("E16").Formula = "=_xlfn.STDEV.S(E10:E14)"
$ws.Range("F16").Formula = "=_xlfn.STDEV.S(F10:F14)"
$ws.Range("G16").Formula = "=_xlfn.STDEV.S(G10:G14)"

# Match the author's final selection in the worksheet.
$ws.Range("G20").Select()
